$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    "('Beast', ['Token Creature — Beast', '4/4'])",
    "('Bird', ['Token Creature — Bird', 'Flying', '1/1'])",
    "('Elemental', ['Token Creature — Elemental', '5/5'])",
    "('Faerie Rogue', ['Token Creature — Faerie Rogue', 'Flying', '1/1'])",
    "('Germ', ['Token Creature — Germ', '0/0'])",
    "('Goblin Rogue', ['Token Creature — Goblin Rogue', '1/1'])",
    "('Kor Ally', ['Token Creature — Kor Ally', '1/1'])",
    "('Rat', ['Token Creature — Rat', '1/1'])",
    "('Saproling', ['Token Creature — Saproling', '1/1'])",
    "('Thopter', ['Token Artifact Creature — Thopter', 'Flying', '1/1'])"
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $values[$i]
}

# Remove the now-unused rows 12 through 37
$ws.Range("A12:A37").EntireRow.Delete()
